$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.669.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.128.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.120.44'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("E9").Value = '  +1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.72%  '
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.648.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.596.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.123.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.00%  '
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0855'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.76%  '
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '446.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0373'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.910.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.280'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.26%  '
$ws.Range("E45").Value = '  +1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.10%  '
